$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.349.13'
$ws.Range('E2').Value = '  +0.11%  '

$ws.Range('D3').Value = '1.935.46'
$ws.Range('E3').Value = '  +0.01%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.43%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7741'
$ws.Range('E5').Value = '  +6.53%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '246.28'
$ws.Range('E6').Value = '  -1.80%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').Value = '  +0.26%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3201'
$ws.Range('E8').Value = '  -3.35%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '27.78'
$ws.Range('E9').Value = '  -0.80%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07055'
$ws.Range('E10').Value = '  -3.28%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7797'
$ws.Range('E11').Value = '  -3.50%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08019'
$ws.Range('E12').Value = '  -0.98%  '

$ws.Range('D13').Value = '1.937.08'
$ws.Range('E13').Value = '  +0.07%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.355'
$ws.Range('E14').Value = '  -2.29%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.66'
$ws.Range('E15').Value = '  -0.18%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.45'
$ws.Range('E16').Value = '  -4.50%  '

$ws.Range('D17').Value = '30.351.68'
$ws.Range('E17').Value = '  +0.11%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '255.61'
$ws.Range('E18').Value = '  +0.83%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007952'
$ws.Range('E19').Value = '  -3.33%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.798'
$ws.Range('E20').Value = '  -0.28%  '

$ws.Range('D21').Value = '2.189.09'
$ws.Range('E21').Value = '  +0.28%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.003'
$ws.Range('E22').Value = '  +0.43%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.004'
$ws.Range('E23').Value = '  +0.48%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.736'
$ws.Range('E24').Value = '  -3.34%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.549'
$ws.Range('E25').Value = '  -2.30%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.69'
$ws.Range('E26').Value = '  -1.01%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1351'
$ws.Range('E27').Value = '  +3.13%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.05'
$ws.Range('E28').Value = '  -1.53%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.273'
$ws.Range('E29').Value = '  -3.68%  '

$ws.Range('E30').Value = '  +1.41%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.515'
$ws.Range('E31').Value = '  -1.67%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.419'
$ws.Range('E32').Value = '  -0.45%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.115'
$ws.Range('E33').Value = '  -2.06%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05152'
$ws.Range('E34').Value = '  -1.93%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.282'
$ws.Range('E35').Value = '  +0.67%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7468'
$ws.Range('E36').Value = '  -0.48%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.779'
$ws.Range('E37').Value = '  +0.71%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01951'
$ws.Range('E38').Value = '  -1.21%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.808'
$ws.Range('E39').Value = '  -0.05%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '78.59'
$ws.Range('E40').Value = '  -0.88%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.428'
$ws.Range('E41').Value = '  -0.24%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4500'
$ws.Range('E42').Value = '  -1.00%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.969'
$ws.Range('E43').Value = '  -3.37%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.002'
$ws.Range('E44').Value = '  +0.30%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8358'
$ws.Range('E45').Value = '  -1.06%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '100.76'
$ws.Range('E46').Value = '  -1.02%  '

$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.817'
$ws.Range('E47').Value = '  +0.67%  '

$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.505'
$ws.Range('E48').Value = '  +0.62%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '985.83'
$ws.Range('E49').Value = '  +10.80%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '37.17'
$ws.Range('E50').Value = '  +0.99%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4146'
$ws.Range('E51').Value = '  -1.17%  '
